# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on "Weekly Quantity" to "Weekly_PO_Qty"
# 2. Rename the "Requested quantity" header on "Monthly Trend" to "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename headers -------------------------------------------------
$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the "PO Forecast" sheet after "Monthly Trend" -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Headers
$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold / border / centered) from an existing header cell
$weekly.Range("A1").Copy()
$forecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = New-Object 'object[,]' 26,4

$data[0,0]=45417.99999999999;  $data[0,1]=8;   $data[0,2]=-99.23977653096017;  $data[0,3]=114.3822611154677
$data[1,0]=45473.99999999999;  $data[1,1]=33;  $data[1,2]=-73.63532576045743;  $data[1,3]=128.8236911141965
$data[2,0]=45487.99999999999;  $data[2,1]=39;  $data[2,2]=-58.74362763247451;  $data[2,3]=149.8293760139039
$data[3,0]=45501.99999999999;  $data[3,1]=45;  $data[3,2]=-61.51433757220882;  $data[3,3]=145.0279194514935
$data[4,0]=45508.99999999999;  $data[4,1]=48;  $data[4,2]=-61.53731000037603;  $data[4,3]=159.8266018352646
$data[5,0]=45515.99999999999;  $data[5,1]=51;  $data[5,2]=-55.34464564743075;  $data[5,3]=154.1148796191494
$data[6,0]=45522.99999999999;  $data[6,1]=54;  $data[6,2]=-52.88041364769571;  $data[6,3]=163.1662985793548
$data[7,0]=45529.99999999999;  $data[7,1]=58;  $data[7,2]=-45.04723841414012;  $data[7,3]=161.5175852419939
$data[8,0]=45543.99999999999;  $data[8,1]=64;  $data[8,2]=-43.52775382604583;  $data[8,3]=167.5104201682381
$data[9,0]=45557.99999999999;  $data[9,1]=70;  $data[9,2]=-35.35971268867063;  $data[9,3]=175.2388085586951
$data[10,0]=45564.99999999999; $data[10,1]=73; $data[10,2]=-32.12496214426726; $data[10,3]=177.4146558537745
$data[11,0]=45571.99999999999; $data[11,1]=76; $data[11,2]=-26.94366184008974; $data[11,3]=183.3015755515291
$data[12,0]=45578.99999999999; $data[12,1]=79; $data[12,2]=-14.3246068715131;  $data[12,3]=181.2509763738042
$data[13,0]=45592.99999999999; $data[13,1]=86; $data[13,2]=-26.45946989766637; $data[13,3]=193.4263326246084
$data[14,0]=45599.99999999999; $data[14,1]=89; $data[14,2]=-12.75619868057782; $data[14,3]=193.4660998549276
$data[15,0]=45627.99999999999; $data[15,1]=101;$data[15,2]=-4.766740259337642; $data[15,3]=199.2060735939936
$data[16,0]=45634.99999999999; $data[16,1]=104;$data[16,2]=3.561830846854331;  $data[16,3]=201.0983572157698
$data[17,0]=45641.99999999999; $data[17,1]=107;$data[17,2]=3.503003538089725;  $data[17,3]=214.3519529706412
$data[18,0]=45648.99999999999; $data[18,1]=110;$data[18,2]=2.985884698650769;  $data[18,3]=220.7789566188835
$data[19,0]=45655.99999999999; $data[19,1]=113;$data[19,2]=11.37536727208895;  $data[19,3]=219.4344863348234
$data[20,0]=45662.99999999999; $data[20,1]=117;$data[20,2]=9.695192984022253;  $data[20,3]=221.4424716745747
$data[21,0]=45669.99999999999; $data[21,1]=120;$data[21,2]=12.58755744316362;  $data[21,3]=216.9924108918464
$data[22,0]=45676.99999999999; $data[22,1]=123;$data[22,2]=18.22621279873009;  $data[22,3]=231.6643247830666
$data[23,0]=45683.99999999999; $data[23,1]=126;$data[23,2]=20.34108048833689;  $data[23,3]=232.7190538892919
$data[24,0]=45690.99999999999; $data[24,1]=129;$data[24,2]=24.92327716017656;  $data[24,3]=235.14532486934
$data[25,0]=45697.99999999999; $data[25,1]=132;$data[25,2]=28.97095034093076;  $data[25,3]=237.6408669126868

$forecast.Range("A2:D27").Value = $data

# Copy the date-number-format style from the "Weekly Quantity" sheet's date
# column onto the new sheet's "ds" column so it matches the workbook's
# existing date-formatted columns.
$weekly.Range("A2").Copy()
$forecast.Range("A2:A27").PasteSpecial(-4122)  # xlPasteFormats

$forecast.Range("A1").Select() | Out-Null
